$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.522.89'
$ws.Range("E2").Value = '  +2.53%  '
$ws.Range("D3").Value = '1.469.90'
$ws.Range("E3").Value = '  +3.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.64%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9545'
$ws.Range("E5").Value = '  -4.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '281.89'
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3700'
$ws.Range("E7").Value = '  -1.38%  '
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.80'
$ws.Range("E9").Value = '  +4.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.056'
$ws.Range("E10").Value = '  +4.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06684'
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.624'
$ws.Range("E13").Value = '  +4.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.27'
$ws.Range("E14").Value = '  +6.95%  '
$ws.Range("D15").Value = '1.474.34'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.267'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001035'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05666'
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.33'
$ws.Range("E19").Value = '  -3.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9545'
$ws.Range("E20").Value = '  -4.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.685'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.71'
$ws.Range("E22").Value = '  +1.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.22'
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.279'
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").Value = '20.675.39'
$ws.Range("E25").Value = '  +3.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.297'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '137.93'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.58'
$ws.Range("E28").Value = '  +4.33%  '
$ws.Range("D29").Value = '1.638.53'
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.91'
$ws.Range("E30").Value = '  +4.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.951'
$ws.Range("E31").Value = '  +2.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.321'
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8358'
$ws.Range("E33").Value = '  -6.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.628'
$ws.Range("E34").Value = '  +27.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07850'
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06043'
$ws.Range("E36").Value = '  +6.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.919'
$ws.Range("E37").Value = '  +3.60%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02074'
$ws.Range("E38").Value = '  +2.93%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.62'
$ws.Range("E39").Value = '  -5.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.122'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9652'
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1889'
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.369'
$ws.Range("E43").Value = '  -12.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5429'
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.55'
$ws.Range("E45").Value = '  +2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.595'
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.40'
$ws.Range("E47").Value = '  +11.25%  '
$ws.Range("E48").Value = '  +4.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.834'
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06444'
$ws.Range("E50").Value = '  +4.37%  '
$ws.Range("E51").Value = '  +0.28%  '
